$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.565.21'
$ws.Range("E2").Value = '  +0.09%  '
$ws.Range("D3").Value = '2.648.77'
$ws.Range("E3").Value = '  -0.08%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.51'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.77'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.95%  '
$ws.Range("E8").Value = '  +0.39%  '
$ws.Range("E9").Value = '  +1.77%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.58'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.68%  '
$ws.Range("E11").Value = '  +4.65%  '
$ws.Range("E12").Value = '  -0.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.58'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.60%  '
$ws.Range("D14").Value = '3.126.94'
$ws.Range("E14").Value = '  -0.04%  '
$ws.Range("D15").Value = '63.423.35'
$ws.Range("E15").Value = '  +0.02%  '
$ws.Range("E16").Value = '  +1.05%  '
$ws.Range("D17").Value = '2.654.78'
$ws.Range("E17").Value = '  -0.05%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.50'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.86%  '
$ws.Range("E19").Value = '  +4.52%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '343.38'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.98%  '
$ws.Range("E21").Value = '  +3.25%  '
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("E23").Value = '  -2.97%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.89'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.40%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.70'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.99%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.06'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +7.67%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '578.50'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +6.74%  '
$ws.Range("E28").Value = '  +2.12%  '
$ws.Range("E29").Value = '  -1.57%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.98'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +2.46%  '
$ws.Range("E31").Value = '  +0.04%  '
$ws.Range("E32").Value = '  +4.35%  '
$ws.Range("E33").Value = '  -3.77%  '
$ws.Range("D34").Value = '0.0₃0823'
$ws.Range("E34").Value = '  +2.08%  '
$ws.Range("E35").Value = '  +7.63%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '168.53'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -3.65%  '
$ws.Range("E37").Value = '  +1.09%  '
$ws.Range("E38").Value = '  -0.14%  '
$ws.Range("E39").Value = '  +7.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.12'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.22%  '
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '169.11'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.00%  '
$ws.Range("E43").Value = '  +1.25%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.18'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0572'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.64%  '
$ws.Range("E46").Value = '  +0.24%  '
$ws.Range("E47").Value = '  +2.96%  '
$ws.Range("E48").Value = '  +0.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.89'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +10.78%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.84'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.31%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.177'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.90%  '
